$d = $word.ActiveDocument

$replacements = @(
    @("14×60=840", "49×21=1029"),
    @("51×45=2295", "54×74=3996"),
    @("65×45=2925", "77×26=2002"),
    @("64×60=3840", "89×62=5518"),
    @("77×19=1463", "17×79=1343"),
    @("45×15=675", "58×34=1972"),
    @("52×47=2444", "97×70=6790"),
    @("96×51=4896", "47×70=3290"),
    @("50×76=3800", "80×42=3360"),
    @("30×81=2430", "65×92=5980"),
    @("48×68=3264", "33×65=2145"),
    @("50×68=3400", "73×50=3650"),
    @("22×43=946", "74×77=5698"),
    @("84×22=1848", "42×46=1932"),
    @("14×50=700", "30×76=2280"),
    @("17×63=1071", "84×36=3024"),
    @("94×80=7520", "94×33=3102"),
    @("53×84=4452", "65×55=3575"),
    @("90×86=7740", "94×61=5734"),
    @("98×38=3724", "28×65=1820"),
    @("22×59=1298", "32×47=1504"),
    @("44×93=4092", "73×72=5256"),
    @("29×34=986", "46×35=1610"),
    @("98×32=3136", "26×58=1508"),
    @("46×87=4002", "80×92=7360")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
